$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.326.70'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.57%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.591.14'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.46%  '
$ws.Range('E4').Value = '  -0.41%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '209.85'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.83%  '
$ws.Range('E6').Value = '  -1.42%  '
$ws.Range('E7').Value = '  -0.38%  '
$ws.Range('E8').Value = '  -1.28%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.245'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.67%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.54'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.90%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0844'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.33%  '
$ws.Range('E12').Value = '  -0.30%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.594.03'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.87%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.08'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.34%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.517'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.60%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.52'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.72%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.350.78'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.32%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0727'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.98%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.46'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '211.12'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.81%  '
$ws.Range('E21').Value = '  -0.47%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.26'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.78%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.16'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.93%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.92'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.17'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.43%  '
$ws.Range('E26').Value = '  -0.43%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.04'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.35%  '
$ws.Range('E28').Value = '  -0.64%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.28'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.39%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0505'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.74%  '
$ws.Range('E31').Value = '  -0.58%  '
$ws.Range('E32').Value = '  -1.68%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.99'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.16%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.305.75'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.44%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.614'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.99%  '
$ws.Range('E36').Value = '  -1.88%  '
$ws.Range('E37').Value = '  -0.88%  '
$ws.Range('E38').Value = '  -0.46%  '
$ws.Range('E39').Value = '  -12.87%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.811'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.72%  '
$ws.Range('E41').Value = '  -0.40%  '
$ws.Range('E42').Value = '  +3.33%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.13'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.75%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '62.61'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.13%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.764'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.74%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.728.07'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.38%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '88.00'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.64%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.50'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.59%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0₆0101'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.61%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0982'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.57%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0504'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.61%  '
